# "Atualizacao de bases das ligas, do dia: 30-05-2024 as 12:21"
#
# The source feed re-sorted a handful of fixtures, which (after the
# re-export) shows up as whole match-records swapping rows, plus a couple
# of standalone HomeTeam/AwayTeam relabels elsewhere on the sheet.
#
# Columns: A=id B=matchId C=Div D=Date E=HomeTeam F=AwayTeam G=FTHG H=FTAG
#          I=ht_goals_h J=ht_goals_a K=FTR L..Q=odds R=Ah S..W=AH odds
#          X..AD=Pinnacle closing values

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Germany Landesliga")

# --- Rows 11 and 12: the two fixtures trade places (id/Div/Date stay put) ---
$ws.Range("B11").Value  = 7035048
$ws.Range("E11").Value  = "SG Unterrath"
$ws.Range("F11").Value  = "TuRU Dsseldorf"
$ws.Range("G11").Value  = 1
$ws.Range("H11").Value  = 0
$ws.Range("I11").Value  = 1
$ws.Range("J11").Value  = 0
$ws.Range("K11").Value  = "H"
$ws.Range("L11").Value  = 3.25
$ws.Range("M11").Value  = 4
$ws.Range("N11").Value  = 1.8
$ws.Range("O11").Value  = 2.9
$ws.Range("P11").Value  = 4
$ws.Range("Q11").Value  = 1.95
$ws.Range("R11").Value  = 0.5
$ws.Range("U11").Value  = 3
$ws.Range("V11").Value  = 1.75
$ws.Range("W11").Value  = 1.95
$ws.Range("X11").Value  = 1.9
$ws.Range("Z11").Value  = -1
$ws.Range("AA11").Value = 0.8
$ws.Range("AB11").Value = -1
$ws.Range("AD11").Value = 0.95

$ws.Range("B12").Value  = 7035046
$ws.Range("E12").Value  = "Cronenberger SC"
$ws.Range("F12").Value  = "FC Viersen"
$ws.Range("G12").Value  = 0
$ws.Range("H12").Value  = 2
$ws.Range("I12").Value  = 0
$ws.Range("J12").Value  = 1
$ws.Range("K12").Value  = "A"
$ws.Range("L12").Value  = 2
$ws.Range("M12").Value  = 3.6
$ws.Range("N12").Value  = 3
$ws.Range("O12").Value  = 2
$ws.Range("P12").Value  = 3.6
$ws.Range("Q12").Value  = 3
$ws.Range("R12").Value  = -0.25
$ws.Range("U12").Value  = 2.75
$ws.Range("V12").Value  = 1.8
$ws.Range("W12").Value  = 2
$ws.Range("X12").Value  = -1
$ws.Range("Z12").Value  = 2
$ws.Range("AA12").Value = -1
$ws.Range("AB12").Value = 1
$ws.Range("AD12").Value = 1

# --- Standalone HomeTeam/AwayTeam relabels ---
$ws.Range("E55").Value = "SG Unterrath"
$ws.Range("E62").Value = "Cronenberger SC"
$ws.Range("F81").Value = "Cronenberger SC"
$ws.Range("E91").Value = "SG Unterrath"
$ws.Range("F112").Value = "SV Rott 1927"
$ws.Range("F116").Value = "BSC Hastedt"

# --- Rows 88, 89, 90: the three fixtures rotate (88<-89, 89<-90, 90<-88) ---
$ws.Range("B88").Value  = 8076438
$ws.Range("E88").Value  = "SV Rott 1927"
$ws.Range("F88").Value  = "SV Breinig"
$ws.Range("G88").Value  = 1
$ws.Range("H88").Value  = 4
$ws.Range("I88").Value  = 1
$ws.Range("J88").Value  = 1
$ws.Range("K88").Value  = "A"
$ws.Range("L88").Value  = 1.363
$ws.Range("M88").Value  = 5
$ws.Range("N88").Value  = 5.5
$ws.Range("O88").Value  = 1.65
$ws.Range("P88").Value  = 4.5
$ws.Range("Q88").Value  = 3.5
$ws.Range("R88").Value  = -0.75
$ws.Range("S88").Value  = 1.825
$ws.Range("T88").Value  = 1.975
$ws.Range("U88").Value  = 3.5
$ws.Range("V88").Value  = 1.975
$ws.Range("W88").Value  = 1.825
$ws.Range("X88").Value  = -1
$ws.Range("Z88").Value  = 2.5
$ws.Range("AB88").Value = 0.9750000000000001
$ws.Range("AC88").Value = 0.9750000000000001

$ws.Range("B89").Value  = 8076477
$ws.Range("E89").Value  = "BSC Hastedt"
$ws.Range("F89").Value  = "SC Weyhe"
$ws.Range("G89").Value  = 2
$ws.Range("H89").Value  = 2
$ws.Range("I89").Value  = 0
$ws.Range("J89").Value  = 0
$ws.Range("K89").Value  = "D"
$ws.Range("L89").Value  = 1.833
$ws.Range("M89").Value  = 4
$ws.Range("N89").Value  = 3.1
$ws.Range("O89").Value  = 1.833
$ws.Range("P89").Value  = 4
$ws.Range("Q89").Value  = 3.1
$ws.Range("R89").Value  = -0.5
$ws.Range("S89").Value  = 1.875
$ws.Range("T89").Value  = 1.925
$ws.Range("U89").Value  = 3.75
$ws.Range("Y89").Value  = 3
$ws.Range("Z89").Value  = -1
$ws.Range("AB89").Value = 0.925
$ws.Range("AC89").Value = 0.4875
$ws.Range("AD89").Value = -0.5

$ws.Range("B90").Value  = 8077795
$ws.Range("E90").Value  = "Eiche Horn"
$ws.Range("F90").Value  = "SVGO Bremen"
$ws.Range("G90").Value  = 5
$ws.Range("H90").Value  = 3
$ws.Range("I90").Value  = 3
$ws.Range("J90").Value  = 2
$ws.Range("K90").Value  = "H"
$ws.Range("L90").Value  = 1.142
$ws.Range("M90").Value  = 7
$ws.Range("N90").Value  = 10
$ws.Range("O90").Value  = 1.083
$ws.Range("P90").Value  = 11
$ws.Range("Q90").Value  = 19
$ws.Range("R90").Value  = -3.5
$ws.Range("S90").Value  = 1.975
$ws.Range("T90").Value  = 1.825
$ws.Range("U90").Value  = 5
$ws.Range("V90").Value  = 1.825
$ws.Range("W90").Value  = 1.975
$ws.Range("X90").Value  = 0.08299999999999996
$ws.Range("Y90").Value  = -1
$ws.Range("AB90").Value = 0.825
$ws.Range("AC90").Value = 0.825
$ws.Range("AD90").Value = -1
